$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.980701754385965
$ws.Range("C2").Value = 0.00701754385964912
$ws.Range("D2").Value = 0.0087719298245614
$ws.Range("E2").Value = 0.0508771929824561
$ws.Range("F2").Value = 0.00526315789473684
$ws.Range("G2").Value = 0.0087719298245614
$ws.Range("H2").Value = 0.975438596491228
$ws.Range("I2").Value = 0.0245614035087719
$ws.Range("J2").Value = 0.954385964912281
$ws.Range("K2").Value = 0.966666666666667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.00350877192982456
$ws.Range("N2").Value = 0.998245614035088
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.996491228070175
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.894736842105263
$ws.Range("S2").Value = 0.0491228070175439
$ws.Range("T2").Value = 0.0614035087719298
$ws.Range("U2").Value = 0.994736842105263
$ws.Range("V2").Value = 0.00526315789473684
$ws.Range("W2").Value = 0.00175438596491228
$ws.Range("X2").Value = 0.0210526315789474
$ws.Range("B3").Value = 0.00526315789473684
$ws.Range("C3").Value = 0.0245614035087719
$ws.Range("D3").Value = 0.00350877192982456
$ws.Range("E3").Value = 0.00526315789473684
$ws.Range("F3").Value = 0.00350877192982456
$ws.Range("G3").Value = 0.978947368421053
$ws.Range("H3").Value = 0.0087719298245614
$ws.Range("I3").Value = 0.012280701754386
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.00175438596491228
$ws.Range("L3").Value = 0.931578947368421
$ws.Range("M3").Value = 0.00526315789473684
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0.998245614035088
$ws.Range("P3").Value = 0.00350877192982456
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 0.087719298245614
$ws.Range("S3").Value = 0.936842105263158
$ws.Range("T3").Value = 0.929824561403509
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0.00526315789473684
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0.00701754385964912
$ws.Range("B4").Value = 0.00350877192982456
$ws.Range("C4").Value = 0.00701754385964912
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.942105263157895
$ws.Range("F4").Value = 0.989473684210526
$ws.Range("G4").Value = 0.0105263157894737
$ws.Range("H4").Value = 0.00526315789473684
$ws.Range("I4").Value = 0.949122807017544
$ws.Range("J4").Value = 0.0087719298245614
$ws.Range("K4").Value = 0.0298245614035088
$ws.Range("L4").Value = 0.00175438596491228
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.00175438596491228
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.0140350877192982
$ws.Range("S4").Value = 0.00350877192982456
$ws.Range("T4").Value = 0.00175438596491228
$ws.Range("U4").Value = 0.00526315789473684
$ws.Range("V4").Value = 0.984210526315789
$ws.Range("W4").Value = 0.996491228070175
$ws.Range("X4").Value = 0.968421052631579
$ws.Range("B5").Value = 0.0105263157894737
$ws.Range("C5").Value = 0.96140350877193
$ws.Range("D5").Value = 0.987719298245614
$ws.Range("E5").Value = 0.00175438596491228
$ws.Range("F5").Value = 0.00175438596491228
$ws.Range("G5").Value = 0.00175438596491228
$ws.Range("H5").Value = 0.0105263157894737
$ws.Range("I5").Value = 0.0140350877192982
$ws.Range("J5").Value = 0.0350877192982456
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.0649122807017544
$ws.Range("M5").Value = 0.991228070175439
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0.00175438596491228
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0.00350877192982456
$ws.Range("S5").Value = 0.0087719298245614
$ws.Range("T5").Value = 0.00701754385964912
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0.00526315789473684
$ws.Range("W5").Value = 0.00175438596491228
$ws.Range("X5").Value = 0.00175438596491228
